$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36, pushing existing rows 36-44 down to 37-45.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with this week's data point.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44468
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108007
$ws.Range("J36").Value = "Coco"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 20
$ws.Range("N36").Value = 24000
$ws.Range("O36").Value = 24000
$ws.Range("P36").Value = 24000
$ws.Range("Q36").Value = "$/malla 20 unidades"
$ws.Range("R36").Value = "Perú"
$ws.Range("S36").Value = 1200
$ws.Range("T36").Value = 20
